$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 3185
$ws.Range("L3").Value = 3267
$ws.Range("K4").Value = 1765
$ws.Range("L4").Value = 825
$ws.Range("L5").Value = 184
$ws.Range("L6").Value = 2888
$ws.Range("K7").Value = 27556
$ws.Range("L7").Value = 10349

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L4").Value = 39
$ws.Range("L7").Value = 347
$ws.Range("L8").Value = 660
$ws.Range("L11").Value = 171
$ws.Range("L15").Value = 77
$ws.Range("K19").Value = 788
$ws.Range("L19").Value = 291
$ws.Range("L20").Value = 261
$ws.Range("L21").Value = 29
$ws.Range("L22").Value = 31
$ws.Range("L23").Value = 110
$ws.Range("L29").Value = 562
$ws.Range("L31").Value = 98
$ws.Range("L33").Value = 476
$ws.Range("L37").Value = 373
$ws.Range("L42").Value = 332
$ws.Range("L48").Value = 140
$ws.Range("L52").Value = 206
$ws.Range("L53").Value = 118
$ws.Range("L54").Value = 217
$ws.Range("L55").Value = 98
$ws.Range("L59").Value = 15
$ws.Range("L63").Value = 33
$ws.Range("L65").Value = 199
$ws.Range("L66").Value = 26
$ws.Range("L67").Value = 379
$ws.Range("L69").Value = 28
$ws.Range("L73").Value = 90
$ws.Range("L75").Value = 39
$ws.Range("L76").Value = 140
$ws.Range("L78").Value = 125
$ws.Range("L79").Value = 266
$ws.Range("L80").Value = 34
$ws.Range("L83").Value = 243
$ws.Range("L85").Value = 526
$ws.Range("L89").Value = 140
$ws.Range("L94").Value = 123
$ws.Range("L95").Value = 142
$ws.Range("L96").Value = 102
$ws.Range("L99").Value = 176
$ws.Range("K101").Value = 27556
$ws.Range("L101").Value = 10349

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L2").Value = 38
$ws.Range("L7").Value = 102

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 108
$ws.Range("L3").Value = 107
$ws.Range("L7").Value = 347

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L3").Value = 53
$ws.Range("L7").Value = 171

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L2").Value = 43
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 140

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 217
$ws.Range("L6").Value = 107
$ws.Range("L7").Value = 526

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 71
$ws.Range("L3").Value = 62
$ws.Range("L7").Value = 206

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("L3").Value = 7
$ws.Range("L7").Value = 28

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L4").Value = 11
$ws.Range("L7").Value = 118

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 194
$ws.Range("L6").Value = 185
$ws.Range("L7").Value = 660

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 78
$ws.Range("L3").Value = 99
$ws.Range("L7").Value = 243

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 134
$ws.Range("L3").Value = 146
$ws.Range("L6").Value = 164
$ws.Range("L7").Value = 476

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L2").Value = 56
$ws.Range("L3").Value = 45
$ws.Range("L6").Value = 28
$ws.Range("L7").Value = 142

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L6").Value = 115
$ws.Range("L7").Value = 373

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 74
$ws.Range("L6").Value = 55
$ws.Range("L7").Value = 199

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 46
$ws.Range("L3").Value = 72
$ws.Range("L6").Value = 41
$ws.Range("L7").Value = 176

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 98

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 144
$ws.Range("L5").Value = 12
$ws.Range("L7").Value = 379

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L2").Value = 46
$ws.Range("L3").Value = 49
$ws.Range("L7").Value = 217

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 211
$ws.Range("L7").Value = 562

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L2").Value = 18
$ws.Range("L6").Value = 58
$ws.Range("L7").Value = 140

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K4").Value = 32
$ws.Range("L6").Value = 91
$ws.Range("K7").Value = 788
$ws.Range("L7").Value = 291

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 27
$ws.Range("L7").Value = 140

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L3").Value = 103
$ws.Range("L6").Value = 93
$ws.Range("L7").Value = 332

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L2").Value = 36
$ws.Range("L7").Value = 125

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L2").Value = 35
$ws.Range("L7").Value = 98

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L2").Value = 32
$ws.Range("L7").Value = 110

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("L2").Value = 4
$ws.Range("L7").Value = 29

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 91
$ws.Range("L3").Value = 93
$ws.Range("L6").Value = 52
$ws.Range("L7").Value = 266

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L6").Value = 73
$ws.Range("L7").Value = 261

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L6").Value = 42
$ws.Range("L7").Value = 123

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L2").Value = 31
$ws.Range("L7").Value = 77

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("L4").Value = 6
$ws.Range("L7").Value = 26

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L4").Value = 7
$ws.Range("L7").Value = 90

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("L3").Value = 7
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("L3").Value = 15
$ws.Range("L7").Value = 39

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("L2").Value = 12
$ws.Range("L7").Value = 31

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("L6").Value = 16
$ws.Range("L7").Value = 34

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("L2").Value = 14
$ws.Range("L7").Value = 39
